$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "0.9993") are stored as literal text, matching the source data,
# instead of being auto-converted to numbers by Excel's input parser.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.306.41"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.930.07"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "0.7549"
$ws.Range("E5").Value = "  +4.54%  "
$ws.Range("D6").Value = "248.31"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("D7").Value = "0.9989"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "28.36"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "0.3211"
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("D10").Value = "0.07099"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").Value = "0.7880"
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("D12").Value = "0.07995"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("D13").Value = "1.929.67"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "5.376"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "94.61"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "14.65"
$ws.Range("E16").Value = "  -2.57%  "
$ws.Range("D17").Value = "30.306.57"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "254.05"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").Value = "0.000008017"
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("D20").Value = "5.810"
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("D21").Value = "2.189.77"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "0.9995"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("D25").Value = "9.577"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("D26").Value = "164.55"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "2.333"
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("D28").Value = "19.11"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").Value = "0.1340"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").Value = "1.360"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").Value = "1.529"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").Value = "4.443"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").Value = "4.142"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").Value = "0.05142"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("D35").Value = "1.290"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "0.7514"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "2.766"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").Value = "0.01967"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").Value = "2.800"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").Value = "78.17"
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("D41").Value = "6.412"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").Value = "0.4515"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").Value = "1.996"
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("D44").Value = "0.9990"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "0.8350"
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("D46").Value = "102.43"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "7.563"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").Value = "9.827"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "989.31"
$ws.Range("E49").Value = "  +12.83%  "
$ws.Range("D50").Value = "37.45"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("D51").Value = "0.1201"
$ws.Range("E51").Value = "  +5.37%  "

# Restore column D's cell style to the original (un-styled) state now
# that the text values are safely stored, so no stray number-format
# style survives into the saved workbook.
$ws.Range("D2:D51").Style = $ws.Range("C2").Style
